$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.996089666666666
$ws.Range("H2").Value = 8.988268999999999
$ws.Range("I2").Value = 0.01916218170515182
$ws.Range("J2").Value = 0.01916218170515182
$ws.Range("M2").Value = 0.04647766666666667
$ws.Range("N2").Value = 0.139433
$ws.Range("O2").Value = 0.1662164546338858
$ws.Range("P2").Value = 0.1662164546338858
$ws.Range("Q2").Value = 0.1392512568307778
$ws.Range("R2").Value = 1.253261311477
$ws.Range("S2").Value = 0.003185069906080644
$ws.Range("T2").Value = 0.003185069906080644
$ws.Range("G3").Value = 2.996089666666666
$ws.Range("H3").Value = 8.988268999999999
$ws.Range("I3").Value = 0.01916218170515182
$ws.Range("J3").Value = 0.01916218170515182
$ws.Range("O3").Value = 0.4511935188540693
$ws.Range("P3").Value = 0.4511935188540693
$ws.Range("Q3").Value = 0.3779966593122221
$ws.Range("R3").Value = 3.401969933809999
$ws.Range("S3").Value = 0.008645852192468518
$ws.Range("T3").Value = 0.008645852192468519
$ws.Range("G4").Value = 2.996089666666666
$ws.Range("H4").Value = 8.988268999999999
$ws.Range("I4").Value = 0.01916218170515182
$ws.Range("J4").Value = 0.01916218170515182
$ws.Range("M4").Value = 0.1069803333333333
$ws.Range("N4").Value = 0.320941
$ws.Range("O4").Value = 0.3825900265120449
$ws.Range("P4").Value = 0.3825900265120449
$ws.Range("Q4").Value = 0.3205226712365555
$ws.Range("R4").Value = 2.884704041129
$ws.Range("S4").Value = 0.007331259606602655
$ws.Range("T4").Value = 0.007331259606602656
$ws.Range("I5").Value = 0.9176106041636097
$ws.Range("J5").Value = 0.9176106041636098
$ws.Range("M5").Value = 0.04647766666666667
$ws.Range("N5").Value = 0.139433
$ws.Range("O5").Value = 0.1662164546338858
$ws.Range("P5").Value = 0.1662164546338858
$ws.Range("Q5").Value = 6.668261050706888
$ws.Range("R5").Value = 60.014349456362
$ws.Range("S5").Value = 0.1525219813585332
$ws.Range("T5").Value = 0.1525219813585332
$ws.Range("I6").Value = 0.9176106041636097
$ws.Range("J6").Value = 0.9176106041636098
$ws.Range("O6").Value = 0.4511935188540693
$ws.Range("P6").Value = 0.4511935188540693
$ws.Range("S6").Value = 0.4140199574303876
$ws.Range("T6").Value = 0.4140199574303876
$ws.Range("I7").Value = 0.9176106041636097
$ws.Range("J7").Value = 0.9176106041636098
$ws.Range("M7").Value = 0.1069803333333333
$ws.Range("N7").Value = 0.320941
$ws.Range("O7").Value = 0.3825900265120449
$ws.Range("P7").Value = 0.3825900265120449
$ws.Range("Q7").Value = 15.34872210936378
$ws.Range("R7").Value = 138.138498984274
$ws.Range("S7").Value = 0.3510686653746889
$ws.Range("T7").Value = 0.351068665374689
$ws.Range("G8").Value = 9.885847333333333
$ws.Range("H8").Value = 29.657542
$ws.Range("I8").Value = 0.06322721413123836
$ws.Range("J8").Value = 0.06322721413123837
$ws.Range("M8").Value = 0.04647766666666667
$ws.Range("N8").Value = 0.139433
$ws.Range("O8").Value = 0.1662164546338858
$ws.Range("P8").Value = 0.1662164546338858
$ws.Range("Q8").Value = 0.4594711170762222
$ws.Range("R8").Value = 4.135240053686
$ws.Range("S8").Value = 0.01050940336927196
$ws.Range("T8").Value = 0.01050940336927197
$ws.Range("G9").Value = 9.885847333333333
$ws.Range("H9").Value = 29.657542
$ws.Range("I9").Value = 0.06322721413123836
$ws.Range("J9").Value = 0.06322721413123837
$ws.Range("O9").Value = 0.4511935188540693
$ws.Range("P9").Value = 0.4511935188540693
$ws.Range("Q9").Value = 1.247231452397778
$ws.Range("R9").Value = 11.22508307158
$ws.Range("S9").Value = 0.02852770923121317
$ws.Range("T9").Value = 0.02852770923121318
$ws.Range("G10").Value = 9.885847333333333
$ws.Range("H10").Value = 29.657542
$ws.Range("I10").Value = 0.06322721413123836
$ws.Range("J10").Value = 0.06322721413123837
$ws.Range("M10").Value = 0.1069803333333333
$ws.Range("N10").Value = 0.320941
$ws.Range("O10").Value = 0.3825900265120449
$ws.Range("P10").Value = 0.3825900265120449
$ws.Range("Q10").Value = 1.057591243002445
$ws.Range("R10").Value = 9.518321187022
$ws.Range("S10").Value = 0.02419010153075322
$ws.Range("T10").Value = 0.02419010153075323
